# Splits each sentence (column B) into an array of words and writes a
# Python-style list-repr string of those words into column C.
# Also removes the extra sentence rows (12-20) that are no longer part
# of the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 12-20 (sentences 10-18) which are dropped in this revision.
$ws.Range("A12:B20").EntireRow.Delete()

# Header for the new column (reuse B1's formatting - bold, centered, bordered).
$ws.Range("C1").Value = "word"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Build the word-list string, e.g. ['foo', 'bar', 'baz'], for each sentence
# in column B (rows 2-11) and place the result into column C.
for ($row = 2; $row -le 11; $row++) {
    $sentence = $ws.Cells.Item($row, 2).Value()
    $words = $sentence.Split(" ")
    $quoted = $words | ForEach-Object { "'" + $_ + "'" }
    $listText = "[" + ($quoted -join ", ") + "]"
    $ws.Cells.Item($row, 3).Value = $listText
}
